$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 135, shifting existing rows 135-214 down to 136-215.
$ws.Rows("135:135").Insert()

# Populate the new row 135 with the new record's data.
$ws.Range("A135").Value = 4
$ws.Range("B135").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C135").Value = "Los Lagos"
$ws.Range("D135").Value = 44603
$ws.Range("E135").Value = 10
$ws.Range("F135").Value = 100112021
$ws.Range("G135").Value = "Ají"
$ws.Range("H135").Value = "Inferno"
$ws.Range("I135").Value = "Primera"
$ws.Range("J135").Value = 160
$ws.Range("K135").Value = 18000
$ws.Range("L135").Value = 18000
$ws.Range("M135").Value = 18000
$ws.Range("N135").Value = "`$/caja 12 kilos"
$ws.Range("O135").Value = "Región de Arica y Parinacota"
$ws.Range("P135").Value = 1500
$ws.Range("Q135").Value = 12
$ws.Range("R135").Value = "Hortaliza"
